# Post pres results commit
# Adds the newly-finished simulation job rows to Sheet1 and refreshes the
# sheet/window chrome (column width, selection, window geometry) to match
# how Excel left things after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New JobID folder rows appended below the existing data (rows 4-9).
$ws.Range("A4").Value = "20240719_MutationStrengthLargerRange"
$ws.Range("A5").Value = "20240722_ExpParam"
$ws.Range("A6").Value = "20240723_ConstantFood"
$ws.Range("A7").Value = "20240723_FracKilledLargerRange"
$ws.Range("A8").Value = "20240723_MetabolicCostLargerRange"
$ws.Range("A9").Value = "20240724_ExpParam_LowerRange"

# Column A was widened to fit the longer folder names.
$ws.Columns("A").ColumnWidth = 34.25

# Resize/reposition the workbook window (as recorded in bookViews).
$win = $excel.ActiveWindow
$win.Left = 11424
$win.Top = 0
$win.Width = 11712
$win.Height = 12336

# Selection left on the first empty row below the new data.
[void]$ws.Range("A10").Select()
